# chore: adapt column header formatting to respective input file names
#
# - Rename the "_old" / "_new" header-name suffixes to the concrete
#   format-version suffixes "_FV2404" / "_FV2410" (columns A:J and L:U of
#   row 1; column K stays "diff").
# - Turn the used range A1:U67 into a real Excel Table ("Table1") so the
#   new headers also drive the table's column names / autofilter.
# - Freeze the header row (row 1) so it stays visible while scrolling.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename header row -------------------------------------------------

$headers = @(
    "Segmentname_FV2404",
    "Segmentgruppe_FV2404",
    "Segment_FV2404",
    "Datenelement_FV2404",
    "Segment ID_FV2404",
    "Code_FV2404",
    "Qualifier_FV2404",
    "Beschreibung_FV2404",
    "Bedingungsausdruck_FV2404",
    "Bedingung_FV2404",
    "diff",
    "Segmentname_FV2410",
    "Segmentgruppe_FV2410",
    "Segment_FV2410",
    "Datenelement_FV2410",
    "Segment ID_FV2410",
    "Code_FV2410",
    "Qualifier_FV2410",
    "Beschreibung_FV2410",
    "Bedingungsausdruck_FV2410",
    "Bedingung_FV2410"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# --- 2. Turn A1:U67 into an Excel Table ------------------------------------

$usedRange = $ws.Range("A1:U67")
$tbl = $ws.ListObjects.Add(1, $usedRange, $null, 1)
$tbl.Name = "Table1"

# --- 3. Freeze the header row ----------------------------------------------

$ws.Cells.Item(2, 1).Select()
$excel.ActiveWindow.FreezePanes = $true

Write-Output "Headers renamed, Table1 created, header row frozen."
